$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update D (Price) / E (Volume 1h) columns with refreshed market data ---
# Note: D-column values that look like plain numbers are written with a leading
# apostrophe so Excel keeps them as text (matching the sheet's original inlineStr
# cell content, e.g. "0.9990" rather than being normalized to 0.999).

$ws.Range("D2").Value2 = "24.813.80"
$ws.Range("E2").Value2 = "  +0.35%  "

$ws.Range("D3").Value2 = "1.709.55"
$ws.Range("E3").Value2 = "  +0.52%  "

$ws.Range("D4").Value2 = "'0.9990"
$ws.Range("E4").Value2 = "  -0.36%  "

$ws.Range("D5").Value2 = "'318.17"
$ws.Range("E5").Value2 = "  +0.72%  "

$ws.Range("D6").Value2 = "'0.9982"
$ws.Range("E6").Value2 = "  -0.35%  "

$ws.Range("D7").Value2 = "'0.3920"
$ws.Range("E7").Value2 = "  -0.36%  "

$ws.Range("D8").Value2 = "'0.4063"
$ws.Range("E8").Value2 = "  +0.32%  "

$ws.Range("E9").Value2 = "  -0.89%  "

$ws.Range("D10").Value2 = "'0.9985"

$ws.Range("D11").Value2 = "'53.46"
$ws.Range("E11").Value2 = "  +1.29%  "

$ws.Range("D12").Value2 = "'0.08836"
$ws.Range("E12").Value2 = "  -0.18%  "

$ws.Range("D13").Value2 = "'26.47"
$ws.Range("E13").Value2 = "  +12.02%  "

$ws.Range("D14").Value2 = "'7.510"
$ws.Range("E14").Value2 = "  -2.30%  "

$ws.Range("D15").Value2 = "'8.135"
$ws.Range("E15").Value2 = "  +0.05%  "

$ws.Range("D16").Value2 = "'0.00001363"
$ws.Range("E16").Value2 = "  +2.88%  "

$ws.Range("D17").Value2 = "1.704.87"

$ws.Range("D18").Value2 = "'97.41"
$ws.Range("E18").Value2 = "  -2.00%  "

$ws.Range("D19").Value2 = "'0.07194"
$ws.Range("E19").Value2 = "  +1.57%  "

$ws.Range("D20").Value2 = "'20.65"
$ws.Range("E20").Value2 = "  +4.22%  "

$ws.Range("D21").Value2 = "'7.321"
$ws.Range("E21").Value2 = "  +2.96%  "

$ws.Range("D22").Value2 = "'0.9984"
$ws.Range("E22").Value2 = "  -0.67%  "

$ws.Range("D23").Value2 = "'14.42"
$ws.Range("E23").Value2 = "  -2.18%  "

$ws.Range("D24").Value2 = "24.819.66"
$ws.Range("E24").Value2 = "  +0.36%  "

$ws.Range("D25").Value2 = "'3.016"
$ws.Range("E25").Value2 = "  -3.78%  "

$ws.Range("D26").Value2 = "'2.330"
$ws.Range("E26").Value2 = "  -0.84%  "

$ws.Range("D27").Value2 = "'23.13"
$ws.Range("E27").Value2 = "  +1.75%  "

$ws.Range("D28").Value2 = "'168.14"
$ws.Range("E28").Value2 = "  +2.09%  "

$ws.Range("D29").Value2 = "'5.978"
$ws.Range("E29").Value2 = "  +16.67%  "

$ws.Range("D30").Value2 = "'8.556"
$ws.Range("E30").Value2 = "  -7.27%  "

$ws.Range("D31").Value2 = "'145.79"
$ws.Range("E31").Value2 = "  +6.96%  "

$ws.Range("D32").Value2 = "1.892.22"
$ws.Range("E32").Value2 = "  +0.74%  "

$ws.Range("D33").Value2 = "'2.195"
$ws.Range("E33").Value2 = "  +11.70%  "

$ws.Range("D34").Value2 = "'0.08842"
$ws.Range("E34").Value2 = "  -2.37%  "

$ws.Range("D35").Value2 = "'1.061"
$ws.Range("E35").Value2 = "  -1.05%  "

$ws.Range("D36").Value2 = "'7.240"
$ws.Range("E36").Value2 = "  -9.40%  "

$ws.Range("D37").Value2 = "'0.03141"
$ws.Range("E37").Value2 = "  +5.01%  "

$ws.Range("E38").Value2 = "  +1.46%  "

# --- Rows 39/40: TheSandbox and FraxShare swapped ranking positions (#38 <-> #39) ---
$ws.Range("B39").Value2 = "TheSandbox"
$ws.Range("C39").Value2 = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value2 = "'0.8486"
$ws.Range("E39").Value2 = "  +8.24%  "

$ws.Range("B40").Value2 = "FraxShare"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value2 = "'10.95"
$ws.Range("E40").Value2 = "  -0.96%  "

$ws.Range("D41").Value2 = "'0.09243"
$ws.Range("E41").Value2 = "  -0.12%  "

$ws.Range("D42").Value2 = "'14.20"
$ws.Range("E42").Value2 = "  -1.49%  "

$ws.Range("D43").Value2 = "'1.481"
$ws.Range("E43").Value2 = "  +1.19%  "

$ws.Range("D44").Value2 = "'17.71"
$ws.Range("E44").Value2 = "  +9.87%  "

$ws.Range("D45").Value2 = "'2.721"
$ws.Range("E45").Value2 = "  +4.12%  "

$ws.Range("D46").Value2 = "'0.7498"
$ws.Range("E46").Value2 = "  +3.21%  "

$ws.Range("D47").Value2 = "'4.283"
$ws.Range("E47").Value2 = "  +1.83%  "

$ws.Range("E48").Value2 = "  +2.81%  "

$ws.Range("D49").Value2 = "'0.9980"
$ws.Range("E49").Value2 = "  -0.36%  "

$ws.Range("D50").Value2 = "'140.66"
$ws.Range("E50").Value2 = "  +0.60%  "

$ws.Range("D51").Value2 = "'0.08298"
$ws.Range("E51").Value2 = "  +3.87%  "
